$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Scaled Stall Speed" (sheet1) : q A/C reference now computed at
# sea-level instead of averaged tunnel dynamic pressure.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Scaled Stall Speed")

$ws1.Range("A34").Value = "q A/C (assume sea level)"
$ws1.Range("B34").Formula = "=0.5*1.225*4.572^2"

# Re-point the Aircraft Lift formulas at the new q (row 34 instead of row 33).
$ws1.Range("I4:I30").Formula = "=`$B`$34*`$K`$1*`$N`$1*H4"
$ws1.Range("I3").Formula = "=`$B`$34*`$K`$1*`$N`$1*H3"

# ---------------------------------------------------------------------------
# Sheet "Scaled 17 fts" (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Scaled 17 fts")

$ws2.Range("B30").Value = "q A/C (assume sea level)"
$ws2.Range("C30").Formula = "=0.5*1.225*5.1816^2"

$ws2.Range("J4:J25").Formula = "=`$C`$30*`$L`$1*`$O`$1*I4"
$ws2.Range("J3").Formula = "=`$C`$30*`$L`$1*`$O`$1*I3"

# ---------------------------------------------------------------------------
# Sheet "Scaled 20 fts" (sheet3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Scaled 20 fts")

$ws3.Range("B36").Value = "q A/C (assume sea level)"
$ws3.Range("C36").Formula = "=0.5*1.225*6.096^2"

$ws3.Range("J4:J32").Formula = "=`$C`$36*`$L`$1*`$O`$1*I4"
$ws3.Range("J3").Formula = "=`$C`$36*`$L`$1*`$O`$1*I3"

# ---------------------------------------------------------------------------
# Selections / active tab: the workbook was left with "Scaled 20 fts" active
# (cell O27 selected), while "Scaled Stall Speed" lost its tab-selected state
# and is left with A34:B34 selected, and "Scaled 17 fts" with B30:C30.
# ---------------------------------------------------------------------------
$null = $ws1.Range("A34:B34").Select()
$null = $ws2.Range("B30:C30").Select()

$null = $ws3.Activate()
$null = $ws3.Range("O27").Select()
